$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write the new row values.
# Columns D and E hold values that LOOK numeric/date-like (small integers
# and dd/mm/yyyy strings) but are stored as plain TEXT in the source sheet
# (same as every pre-existing data row). A leading apostrophe forces text
# entry so Excel does not reinterpret "03/01/2026" as a real date serial
# (which would flip it to the US mm/dd reading) or "1".."10" as numbers.
$ws.Cells.Item(244,1).Value = 'Year 5'
$ws.Cells.Item(244,2).Value = 'B2-D1'
$ws.Cells.Item(244,3).Value = 'endocrinology'
$ws.Cells.Item(244,4).Value = "'1"
$ws.Cells.Item(244,5).Value = "'20/12/2025"
$ws.Cells.Item(244,6).Value = '09:00:00'
$ws.Cells.Item(244,7).Value = 360

$ws.Cells.Item(245,1).Value = 'Year 5'
$ws.Cells.Item(245,2).Value = 'B2-D1'
$ws.Cells.Item(245,3).Value = 'endocrinology'
$ws.Cells.Item(245,4).Value = "'2"
$ws.Cells.Item(245,5).Value = "'21/12/2025"
$ws.Cells.Item(245,6).Value = '09:00:00'
$ws.Cells.Item(245,7).Value = 360

$ws.Cells.Item(246,1).Value = 'Year 5'
$ws.Cells.Item(246,2).Value = 'B2-D1'
$ws.Cells.Item(246,3).Value = 'endocrinology'
$ws.Cells.Item(246,4).Value = "'3"
$ws.Cells.Item(246,5).Value = "'22/12/2025"
$ws.Cells.Item(246,6).Value = '09:00:00'
$ws.Cells.Item(246,7).Value = 360

$ws.Cells.Item(247,1).Value = 'Year 5'
$ws.Cells.Item(247,2).Value = 'B2-D1'
$ws.Cells.Item(247,3).Value = 'endocrinology'
$ws.Cells.Item(247,4).Value = "'4"
$ws.Cells.Item(247,5).Value = "'23/12/2025"
$ws.Cells.Item(247,6).Value = '09:00:00'
$ws.Cells.Item(247,7).Value = 360

$ws.Cells.Item(248,1).Value = 'Year 5'
$ws.Cells.Item(248,2).Value = 'B2-D1'
$ws.Cells.Item(248,3).Value = 'endocrinology'
$ws.Cells.Item(248,4).Value = "'5"
$ws.Cells.Item(248,5).Value = "'24/12/2025"
$ws.Cells.Item(248,6).Value = '09:00:00'
$ws.Cells.Item(248,7).Value = 360

$ws.Cells.Item(249,1).Value = 'Year 5'
$ws.Cells.Item(249,2).Value = 'B2-D1'
$ws.Cells.Item(249,3).Value = 'endocrinology'
$ws.Cells.Item(249,4).Value = "'6"
$ws.Cells.Item(249,5).Value = "'27/12/2025"
$ws.Cells.Item(249,6).Value = '09:00:00'
$ws.Cells.Item(249,7).Value = 360

$ws.Cells.Item(250,1).Value = 'Year 5'
$ws.Cells.Item(250,2).Value = 'B2-D1'
$ws.Cells.Item(250,3).Value = 'endocrinology'
$ws.Cells.Item(250,4).Value = "'7"
$ws.Cells.Item(250,5).Value = "'28/12/2025"
$ws.Cells.Item(250,6).Value = '09:00:00'
$ws.Cells.Item(250,7).Value = 360

$ws.Cells.Item(251,1).Value = 'Year 5'
$ws.Cells.Item(251,2).Value = 'B2-D1'
$ws.Cells.Item(251,3).Value = 'endocrinology'
$ws.Cells.Item(251,4).Value = "'8"
$ws.Cells.Item(251,5).Value = "'29/12/2025"
$ws.Cells.Item(251,6).Value = '09:00:00'
$ws.Cells.Item(251,7).Value = 360

$ws.Cells.Item(252,1).Value = 'Year 5'
$ws.Cells.Item(252,2).Value = 'B2-D1'
$ws.Cells.Item(252,3).Value = 'endocrinology'
$ws.Cells.Item(252,4).Value = "'9"
$ws.Cells.Item(252,5).Value = "'30/12/2025"
$ws.Cells.Item(252,6).Value = '09:00:00'
$ws.Cells.Item(252,7).Value = 360

$ws.Cells.Item(253,1).Value = 'Year 5'
$ws.Cells.Item(253,2).Value = 'B2-D1'
$ws.Cells.Item(253,3).Value = 'endocrinology'
$ws.Cells.Item(253,4).Value = "'10"
$ws.Cells.Item(253,5).Value = "'31/12/2025"
$ws.Cells.Item(253,6).Value = '09:00:00'
$ws.Cells.Item(253,7).Value = 360

$ws.Cells.Item(254,1).Value = 'Year 5'
$ws.Cells.Item(254,2).Value = 'B2-D1'
$ws.Cells.Item(254,3).Value = 'gastroenterology'
$ws.Cells.Item(254,4).Value = "'1"
$ws.Cells.Item(254,5).Value = "'03/01/2026"
$ws.Cells.Item(254,6).Value = '09:00:00'
$ws.Cells.Item(254,7).Value = 360

$ws.Cells.Item(255,1).Value = 'Year 5'
$ws.Cells.Item(255,2).Value = 'B2-D1'
$ws.Cells.Item(255,3).Value = 'gastroenterology'
$ws.Cells.Item(255,4).Value = "'2"
$ws.Cells.Item(255,5).Value = "'04/01/2026"
$ws.Cells.Item(255,6).Value = '09:00:00'
$ws.Cells.Item(255,7).Value = 360

$ws.Cells.Item(256,1).Value = 'Year 5'
$ws.Cells.Item(256,2).Value = 'B2-D1'
$ws.Cells.Item(256,3).Value = 'gastroenterology'
$ws.Cells.Item(256,4).Value = "'3"
$ws.Cells.Item(256,5).Value = "'05/01/2026"
$ws.Cells.Item(256,6).Value = '09:00:00'
$ws.Cells.Item(256,7).Value = 360

$ws.Cells.Item(257,1).Value = 'Year 5'
$ws.Cells.Item(257,2).Value = 'B2-D1'
$ws.Cells.Item(257,3).Value = 'gastroenterology'
$ws.Cells.Item(257,4).Value = "'4"
$ws.Cells.Item(257,5).Value = "'06/01/2026"
$ws.Cells.Item(257,6).Value = '09:00:00'
$ws.Cells.Item(257,7).Value = 360

$ws.Cells.Item(258,1).Value = 'Year 5'
$ws.Cells.Item(258,2).Value = 'B2-D1'
$ws.Cells.Item(258,3).Value = 'gastroenterology'
$ws.Cells.Item(258,4).Value = "'5"
$ws.Cells.Item(258,5).Value = "'07/01/2026"
$ws.Cells.Item(258,6).Value = '09:00:00'
$ws.Cells.Item(258,7).Value = 360

$ws.Cells.Item(259,1).Value = 'Year 5'
$ws.Cells.Item(259,2).Value = 'B2-D1'
$ws.Cells.Item(259,3).Value = 'gastroenterology'
$ws.Cells.Item(259,4).Value = "'6"
$ws.Cells.Item(259,5).Value = "'10/01/2026"
$ws.Cells.Item(259,6).Value = '09:00:00'
$ws.Cells.Item(259,7).Value = 360

$ws.Cells.Item(260,1).Value = 'Year 5'
$ws.Cells.Item(260,2).Value = 'B2-D1'
$ws.Cells.Item(260,3).Value = 'gastroenterology'
$ws.Cells.Item(260,4).Value = "'7"
$ws.Cells.Item(260,5).Value = "'11/01/2026"
$ws.Cells.Item(260,6).Value = '09:00:00'
$ws.Cells.Item(260,7).Value = 360

$ws.Cells.Item(261,1).Value = 'Year 5'
$ws.Cells.Item(261,2).Value = 'B2-D1'
$ws.Cells.Item(261,3).Value = 'gastroenterology'
$ws.Cells.Item(261,4).Value = "'8"
$ws.Cells.Item(261,5).Value = "'12/01/2026"
$ws.Cells.Item(261,6).Value = '09:00:00'
$ws.Cells.Item(261,7).Value = 360

$ws.Cells.Item(262,1).Value = 'Year 5'
$ws.Cells.Item(262,2).Value = 'B2-D1'
$ws.Cells.Item(262,3).Value = 'gastroenterology'
$ws.Cells.Item(262,4).Value = "'9"
$ws.Cells.Item(262,5).Value = "'13/01/2026"
$ws.Cells.Item(262,6).Value = '09:00:00'
$ws.Cells.Item(262,7).Value = 360

$ws.Cells.Item(263,1).Value = 'Year 5'
$ws.Cells.Item(263,2).Value = 'B2-D1'
$ws.Cells.Item(263,3).Value = 'gastroenterology'
$ws.Cells.Item(263,4).Value = "'10"
$ws.Cells.Item(263,5).Value = "'14/01/2026"
$ws.Cells.Item(263,6).Value = '09:00:00'
$ws.Cells.Item(263,7).Value = 360

$ws.Cells.Item(264,1).Value = 'Year 5'
$ws.Cells.Item(264,2).Value = 'B2-D1'
$ws.Cells.Item(264,3).Value = 'nephrology'
$ws.Cells.Item(264,4).Value = "'1"
$ws.Cells.Item(264,5).Value = "'07/02/2026"
$ws.Cells.Item(264,6).Value = '09:00:00'
$ws.Cells.Item(264,7).Value = 360

$ws.Cells.Item(265,1).Value = 'Year 5'
$ws.Cells.Item(265,2).Value = 'B2-D1'
$ws.Cells.Item(265,3).Value = 'nephrology'
$ws.Cells.Item(265,4).Value = "'2"
$ws.Cells.Item(265,5).Value = "'08/02/2026"
$ws.Cells.Item(265,6).Value = '09:00:00'
$ws.Cells.Item(265,7).Value = 360

$ws.Cells.Item(266,1).Value = 'Year 5'
$ws.Cells.Item(266,2).Value = 'B2-D1'
$ws.Cells.Item(266,3).Value = 'nephrology'
$ws.Cells.Item(266,4).Value = "'3"
$ws.Cells.Item(266,5).Value = "'09/02/2026"
$ws.Cells.Item(266,6).Value = '09:00:00'
$ws.Cells.Item(266,7).Value = 360

$ws.Cells.Item(267,1).Value = 'Year 5'
$ws.Cells.Item(267,2).Value = 'B2-D1'
$ws.Cells.Item(267,3).Value = 'nephrology'
$ws.Cells.Item(267,4).Value = "'4"
$ws.Cells.Item(267,5).Value = "'10/02/2026"
$ws.Cells.Item(267,6).Value = '09:00:00'
$ws.Cells.Item(267,7).Value = 360

$ws.Cells.Item(268,1).Value = 'Year 5'
$ws.Cells.Item(268,2).Value = 'B2-D1'
$ws.Cells.Item(268,3).Value = 'nephrology'
$ws.Cells.Item(268,4).Value = "'5"
$ws.Cells.Item(268,5).Value = "'11/02/2026"
$ws.Cells.Item(268,6).Value = '09:00:00'
$ws.Cells.Item(268,7).Value = 360

$ws.Cells.Item(269,1).Value = 'Year 5'
$ws.Cells.Item(269,2).Value = 'B2-D1'
$ws.Cells.Item(269,3).Value = 'neurology'
$ws.Cells.Item(269,4).Value = "'1"
$ws.Cells.Item(269,5).Value = "'06/12/2025"
$ws.Cells.Item(269,6).Value = '09:00:00'
$ws.Cells.Item(269,7).Value = 360

$ws.Cells.Item(270,1).Value = 'Year 5'
$ws.Cells.Item(270,2).Value = 'B2-D1'
$ws.Cells.Item(270,3).Value = 'neurology'
$ws.Cells.Item(270,4).Value = "'2"
$ws.Cells.Item(270,5).Value = "'07/12/2025"
$ws.Cells.Item(270,6).Value = '09:00:00'
$ws.Cells.Item(270,7).Value = 360

$ws.Cells.Item(271,1).Value = 'Year 5'
$ws.Cells.Item(271,2).Value = 'B2-D1'
$ws.Cells.Item(271,3).Value = 'neurology'
$ws.Cells.Item(271,4).Value = "'3"
$ws.Cells.Item(271,5).Value = "'08/12/2025"
$ws.Cells.Item(271,6).Value = '09:00:00'
$ws.Cells.Item(271,7).Value = 360

$ws.Cells.Item(272,1).Value = 'Year 5'
$ws.Cells.Item(272,2).Value = 'B2-D1'
$ws.Cells.Item(272,3).Value = 'neurology'
$ws.Cells.Item(272,4).Value = "'4"
$ws.Cells.Item(272,5).Value = "'09/12/2025"
$ws.Cells.Item(272,6).Value = '09:00:00'
$ws.Cells.Item(272,7).Value = 360

$ws.Cells.Item(273,1).Value = 'Year 5'
$ws.Cells.Item(273,2).Value = 'B2-D1'
$ws.Cells.Item(273,3).Value = 'neurology'
$ws.Cells.Item(273,4).Value = "'5"
$ws.Cells.Item(273,5).Value = "'13/12/2025"
$ws.Cells.Item(273,6).Value = '09:00:00'
$ws.Cells.Item(273,7).Value = 360

$ws.Cells.Item(274,1).Value = 'Year 5'
$ws.Cells.Item(274,2).Value = 'B2-D1'
$ws.Cells.Item(274,3).Value = 'neurology'
$ws.Cells.Item(274,4).Value = "'6"
$ws.Cells.Item(274,5).Value = "'14/12/2025"
$ws.Cells.Item(274,6).Value = '09:00:00'
$ws.Cells.Item(274,7).Value = 360

$ws.Cells.Item(275,1).Value = 'Year 5'
$ws.Cells.Item(275,2).Value = 'B2-D1'
$ws.Cells.Item(275,3).Value = 'neurology'
$ws.Cells.Item(275,4).Value = "'7"
$ws.Cells.Item(275,5).Value = "'15/12/2025"
$ws.Cells.Item(275,6).Value = '09:00:00'
$ws.Cells.Item(275,7).Value = 360

$ws.Cells.Item(276,1).Value = 'Year 5'
$ws.Cells.Item(276,2).Value = 'B2-D1'
$ws.Cells.Item(276,3).Value = 'neurology'
$ws.Cells.Item(276,4).Value = "'8"
$ws.Cells.Item(276,5).Value = "'16/12/2025"
$ws.Cells.Item(276,6).Value = '09:00:00'
$ws.Cells.Item(276,7).Value = 360

$ws.Cells.Item(277,1).Value = 'Year 5'
$ws.Cells.Item(277,2).Value = 'B2-D1'
$ws.Cells.Item(277,3).Value = 'physical medicine'
$ws.Cells.Item(277,4).Value = "'1"
$ws.Cells.Item(277,5).Value = "'10/12/2025"
$ws.Cells.Item(277,6).Value = '09:00:00'
$ws.Cells.Item(277,7).Value = 360

$ws.Cells.Item(278,1).Value = 'Year 5'
$ws.Cells.Item(278,2).Value = 'B2-D1'
$ws.Cells.Item(278,3).Value = 'physical medicine'
$ws.Cells.Item(278,4).Value = "'2"
$ws.Cells.Item(278,5).Value = "'17/12/2025"
$ws.Cells.Item(278,6).Value = '09:00:00'
$ws.Cells.Item(278,7).Value = 360

$ws.Cells.Item(279,1).Value = 'Year 5'
$ws.Cells.Item(279,2).Value = 'B2-D1'
$ws.Cells.Item(279,3).Value = 'rheumatology'
$ws.Cells.Item(279,4).Value = "'1"
$ws.Cells.Item(279,5).Value = "'17/01/2026"
$ws.Cells.Item(279,6).Value = '09:00:00'
$ws.Cells.Item(279,7).Value = 360

$ws.Cells.Item(280,1).Value = 'Year 5'
$ws.Cells.Item(280,2).Value = 'B2-D1'
$ws.Cells.Item(280,3).Value = 'rheumatology'
$ws.Cells.Item(280,4).Value = "'2"
$ws.Cells.Item(280,5).Value = "'18/01/2026"
$ws.Cells.Item(280,6).Value = '09:00:00'
$ws.Cells.Item(280,7).Value = 360

$ws.Cells.Item(281,1).Value = 'Year 5'
$ws.Cells.Item(281,2).Value = 'B2-D1'
$ws.Cells.Item(281,3).Value = 'rheumatology'
$ws.Cells.Item(281,4).Value = "'3"
$ws.Cells.Item(281,5).Value = "'19/01/2026"
$ws.Cells.Item(281,6).Value = '09:00:00'
$ws.Cells.Item(281,7).Value = 360

$ws.Cells.Item(282,1).Value = 'Year 5'
$ws.Cells.Item(282,2).Value = 'B2-D1'
$ws.Cells.Item(282,3).Value = 'rheumatology'
$ws.Cells.Item(282,4).Value = "'4"
$ws.Cells.Item(282,5).Value = "'20/01/2026"
$ws.Cells.Item(282,6).Value = '09:00:00'
$ws.Cells.Item(282,7).Value = 360

$ws.Cells.Item(283,1).Value = 'Year 5'
$ws.Cells.Item(283,2).Value = 'B2-D1'
$ws.Cells.Item(283,3).Value = 'rheumatology'
$ws.Cells.Item(283,4).Value = "'5"
$ws.Cells.Item(283,5).Value = "'21/01/2026"
$ws.Cells.Item(283,6).Value = '09:00:00'
$ws.Cells.Item(283,7).Value = 360

# Step 2: copy the alternating-row formatting (styles 2,3,4,5 / 6,7,8,9)
# down from the last two existing data rows onto the newly added rows.
$src = $ws.Range("A242:G243")
$dst = $ws.Range("A244:G283")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Added rows 244-283"